$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns used in this block (skips G/H, which are unaffected)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row (2..25), one array per row aligned to $cols
$data = @{
    2 = @(1.02, 1.039051041329276, 1.044324510941993, 1.047167603258968, 1.056942686806193, 1.043462659208072, 1.044145154734633, 1.047095880588599, 1.049930997312892, 1.059679013980852, 1.04562796061948)
    3 = @(1.02, 1.039995381693004, 1.044911987269564, 1.048017288673664, 1.057899958650401, 1.043701229512901, 1.044734566968081, 1.047495155965903, 1.050592374062448, 1.060449659374841, 1.046218209885906)
    4 = @(1.02, 1.040606827416793, 1.045292017208368, 1.048567824946624, 1.058520220700019, 1.043854303088766, 1.045115726798075, 1.047752708253473, 1.05102041693274, 1.060948537097056, 1.046599911006578)
    5 = @(1.02, 1.040863972663084, 1.045451754619064, 1.048799444566519, 1.05878117938034, 1.043918343901271, 1.045275910565849, 1.047860789186525, 1.051200386156108, 1.06115831645893, 1.046760322253675)
    6 = @(1.02, 1.040907153934751, 1.045478573590878, 1.048838344695273, 1.058825007210428, 1.043929078358453, 1.045302802824674, 1.047878925036971, 1.051230604923241, 1.061193542303325, 1.046787252702592)
    7 = @(1.02, 1.040610263038088, 1.045294151736361, 1.048570919176592, 1.058523706856469, 1.043855160029418, 1.04511786740194, 1.04775415320017, 1.051022821613026, 1.060951339979571, 1.046602054650345)
    8 = @(1.02, 1.039370103351027, 1.04452307237101, 1.047454605865056, 1.057266025922091, 1.043543553350365, 1.044344396276844, 1.047230983657317, 1.050154493858884, 1.059939411229087, 1.045827485107534)
    9 = @(1.02, 1.037187845622992, 1.043163589810153, 1.045493190770747, 1.055056348027313, 1.042984558352082, 1.042979726177853, 1.04630297476448, 1.048625108378085, 1.058157989475196, 1.044460877020472)
    10 = @(1.02, 1.035735120756845, 1.042256868086695, 1.044189470963001, 1.053587689392306, 1.042605278724324, 1.042068847376143, 1.045680259555176, 1.047606062451701, 1.056971604108349, 1.043548704666384)
    11 = @(1.02, 1.035106587266124, 1.041864173470877, 1.043625884581509, 1.052952818074108, 1.04243948601177, 1.041674177398867, 1.045409672834149, 1.047164946040878, 1.056458191864999, 1.043153474212472)
    12 = @(1.02, 1.034873198792367, 1.041718299132186, 1.043416684923062, 1.052717160262294, 1.042377669202586, 1.041527542184363, 1.045309023772943, 1.047001117360228, 1.056267533943309, 1.043006630759144)
    13 = @(1.02, 1.034923257921168, 1.041749590102216, 1.043461552540286, 1.052767702326811, 1.042390939696721, 1.041558997605679, 1.045330619704948, 1.047036258175361, 1.056308428604259, 1.043038130850766)
    14 = @(1.02, 1.035087293729513, 1.041852115648433, 1.043608589180326, 1.052933335220613, 1.042434380984925, 1.041662057240407, 1.045401356023158, 1.047151403464815, 1.05644243107193, 1.043141336841997)
    15 = @(1.02, 1.035188371813338, 1.041915283736043, 1.043699202083927, 1.053035408573164, 1.042461115607531, 1.041725550777662, 1.045444920340242, 1.047222351168865, 1.056525000575933, 1.043204920547357)
    16 = @(1.02, 1.035776845185836, 1.042282928431173, 1.044226894097657, 1.053629846311855, 1.042616248986261, 1.042095035056743, 1.045698197653116, 1.047635340869589, 1.057005684012576, 1.043574929536495)
    17 = @(1.02, 1.036146115034075, 1.042513522495256, 1.04455815200935, 1.054003008116579, 1.042713142308529, 1.042326735539403, 1.045856819120678, 1.047894435823512, 1.057307285252388, 1.043806959060417)
    18 = @(1.02, 1.036361552650892, 1.042648016462189, 1.044751459071509, 1.054220770002679, 1.042769507775618, 1.042461857999051, 1.045949248896794, 1.048045574698515, 1.057483233067693, 1.043942273409456)
    19 = @(1.02, 1.036435019556334, 1.042693874081116, 1.044817386998998, 1.054295038596427, 1.042788701349681, 1.042507927081973, 1.045980749516329, 1.048097111353419, 1.057543231596626, 1.04398840791576)
    20 = @(1.02, 1.036106490824345, 1.042488782713673, 1.044522601859751, 1.053962960731889, 1.042702762161814, 1.042301878799502, 1.045839809988298, 1.047866636017451, 1.057274923294868, 1.043782067021094)
    21 = @(1.02, 1.035038987144276, 1.041821924695202, 1.0435652866407, 1.052884555988741, 1.042421595062772, 1.041631709766966, 1.045380529826395, 1.047117495442973, 1.056402969399646, 1.043110946271661)
    22 = @(1.02, 1.034368250259032, 1.041402587127963, 1.042964203387352, 1.0522074560911, 1.042243460552288, 1.041210133027281, 1.045090946340998, 1.046646606459308, 1.055855005486335, 1.042688770844651)
    23 = @(1.02, 1.034723778107477, 1.041624890799266, 1.043282770957446, 1.052566310413493, 1.042338021105343, 1.041433638911681, 1.045244536969406, 1.046896221401758, 1.056145465741526, 1.042912594133043)
    24 = @(1.02, 1.036124395148485, 1.042499961579888, 1.044538665164376, 1.053981056103319, 1.042707452969214, 1.04231311056612, 1.045847495965736, 1.047879197512564, 1.057289546181216, 1.043793314738109)
    25 = @(1.02, 1.037751643041339, 1.043515126371304, 1.045999583539414, 1.055626822384501, 1.043130241235352, 1.043332723755322, 1.046543604983117, 1.049020399697364, 1.058618318084539, 1.044814375894991)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}